$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.338.14"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "3.373.92"
$ws.Range("E3").Value = "  -2.18%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'567.23"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("D6").Value = "'140.17"
$ws.Range("E6").Value = "  -6.47%  "

$ws.Range("D8").Value = "3.374.77"
$ws.Range("E8").Value = "  -2.18%  "

$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").Value = "'7.48"
$ws.Range("E10").Value = "  -3.90%  "

$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("D13").Value = "3.950.95"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'28.07"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.124"
$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").Value = "3.371.58"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("E17").Value = "  -3.72%  "

$ws.Range("D18").Value = "60.478.03"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("E20").Value = "  -4.12%  "

$ws.Range("D21").Value = "'8.99"
$ws.Range("E21").Value = "  -5.60%  "

$ws.Range("D22").Value = "'385.81"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").Value = "'0.554"
$ws.Range("E23").Value = "  -2.47%  "

$ws.Range("D24").Value = "'73.09"

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  -8.02%  "

$ws.Range("D27").Value = "3.519.05"
$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'7.34"
$ws.Range("E30").Value = "  -5.06%  "

$ws.Range("D31").Value = "'7.92"
$ws.Range("E31").Value = "  -4.08%  "

$ws.Range("E32").Value = "  -2.36%  "

$ws.Range("E33").Value = "  -9.03%  "

$ws.Range("D35").Value = "'23.49"
$ws.Range("E35").Value = "  -2.50%  "

$ws.Range("D36").Value = "3.404.27"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").Value = "'6.89"
$ws.Range("E37").Value = "  -2.39%  "

$ws.Range("D38").Value = "'168.43"
$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").Value = "'4.91"
$ws.Range("E39").Value = "  -5.83%  "

$ws.Range("E40").Value = "  -4.92%  "

$ws.Range("D41").Value = "'0.0769"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").Value = "'27.17"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'0.777"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("E45").Value = "  -1.94%  "

$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").Value = "'41.26"
$ws.Range("E47").Value = "  -2.64%  "

$ws.Range("D48").Value = "2.515.35"
$ws.Range("E48").Value = "  -3.13%  "

$ws.Range("E49").Value = "  -4.68%  "

$ws.Range("D50").Value = "'23.09"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("E51").Value = "  -3.30%  "
